# Repull data, push all data, mean calculation
# Update column F ("dSF") values for a set of rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    3  = 8
    5  = -5
    15 = 4
    19 = 1
    24 = -1
    26 = -5
    31 = -8
    33 = 4
    38 = -3
    39 = -6
    40 = -7
    43 = -3
    46 = -4
    55 = -2
    57 = -4
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
